$d = $word.ActiveDocument

# The document has a short "M" bullet (numId=7 list) right before a
# "_GoBack" bookmark. Expand it to "Java Android" and add a new sibling
# bullet "Mobile" directly after it, keeping the bookmark anchored to the
# very end of the list (now trailing the new "Mobile" run), exactly as in
# the target revision.

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "M") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $r = $target.Range

    # Replace "M" with "Java Android" (keep the trailing paragraph mark / bookmark).
    $textRange = $d.Range($r.Start, $r.Start + 1)
    $textRange.Text = "Java Android"

    # Re-fetch the (now longer) paragraph range and split it right before the
    # paragraph mark, so the bookmark that sits there ends up in the new,
    # second paragraph.
    $p2 = $target
    $r2 = $p2.Range
    $splitPoint = $d.Range($r2.End - 1, $r2.End - 1)
    $splitPoint.InsertBefore([char]13)

    # The new (second) paragraph is now empty except for the relocated
    # bookmark; insert the "Mobile" run before that bookmark so the bookmark
    # still trails the text, matching the source paragraph's original layout.
    $newPara = $p2.Next()
    $newPara.Range.InsertBefore("Mobile")
}
